$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value of B5 to "UML" (new entry in the Team Member / EDT table)
$ws.Range("B5").Value = "UML"

# Move the selection to B5, as in the saved worksheet view
$ws.Range("B5").Select() | Out-Null
